# Refresh the "token"/"id" columns (C2:D4) on the one worksheet with the
# latest batch of login credentials generated by the test run.
#
# Each row corresponds to one test user (daniel5f, Jorge2525, mario35):
#   column C = freshly minted JWT auth token
#   column D = freshly minted UUID session/record id

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - daniel5f
$ws.Range("C2").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6ImRhbmllbDVmIiwicGFzc3dvcmQiOiJBejI1Mjg4QCIsImlhdCI6MTcwMjE0MjIwOX0.dunWqUFAeSD_q61JQ9i1JljDIHawqHpOznMSnlJWWqE"
$ws.Range("D2").Value = "7817e48d-65c8-4992-9427-411433621c4a"

# Row 3 - Jorge2525
$ws.Range("C3").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6IkpvcmdlMjUyNSIsInBhc3N3b3JkIjoiYXNUMzU2NDQ0QCIsImlhdCI6MTcwMjE0MjIxMH0.AwlogkkBFOAKNJEyto4hNFreNJskJTOsNvygH7NfKJE"
$ws.Range("D3").Value = "1f7a37d5-73c7-49a2-9d77-c649bb24a408"

# Row 4 - mario35
$ws.Range("C4").Value = "eyJhbGciOiJIUzI1NiIsInR5cCI6IkpXVCJ9.eyJ1c2VyTmFtZSI6Im1hcmlvMzUiLCJwYXNzd29yZCI6Im1BcmlvdXVnQDMiLCJpYXQiOjE3MDIxNDIyMTJ9.-rtLLD7NDyNvzIczE3Gk_g7-HhbguAhGk85edNNptl4"
$ws.Range("D4").Value = "2f6c52dc-b44c-4d04-a2c5-157b08a00bcc"
